$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 730
$ws.Range("F4").Value = 1436
$ws.Range("F5").Value = 217
$ws.Range("F6").Value = 86
$ws.Range("F7").Value = 127
$ws.Range("F8").Value = 6104
$ws.Range("F9").Value = 67
$ws.Range("F12").Value = 4941
$ws.Range("F15").Value = 1163
$ws.Range("F17").Value = 351
$ws.Range("F18").Value = 55
$ws.Range("F20").Value = 283
$ws.Range("F21").Value = 22
$ws.Range("F22").Value = 3440
$ws.Range("F23").Value = 141

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 64

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 64
$ws.Range("F4").Value = 730
$ws.Range("F5").Value = 1436
$ws.Range("F6").Value = 217
$ws.Range("F7").Value = 86
$ws.Range("F8").Value = 127
$ws.Range("F9").Value = 6104
$ws.Range("F10").Value = 67
$ws.Range("F13").Value = 4941
$ws.Range("F16").Value = 1163
$ws.Range("F18").Value = 351
$ws.Range("F19").Value = 55
$ws.Range("F21").Value = 283
$ws.Range("F22").Value = 22
$ws.Range("F23").Value = 3440
$ws.Range("F25").Value = 141
